{"js": "// The edit:\n//  1. Three Heading2 paragraphs (\"Introduction\", \"Immediate Effects of\n//     Social Media\", \"Conclusion\") lose their heading style and revert to\n//     the document's default (\"Normal\") paragraph style.\n//  2. Four in-text citations are swapped out for new ones (a citation\n//     check / reference update):\n//       (Karim and Oyewande)  -> (Lee 208)\n//       (Huang)                -> (Lee 208)\n//       (Jiang and Ngien)      -> (Ref-f422876)\n//       (Lee)                  -> (Ref-f422876)\n//       (Naslund et al.)       -> (Brown & Garcia, 2018)   [occurs twice]\n\n// --- 1. Demote the three Heading 2 paragraphs back to Normal -------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].style === \"Heading 2\") {\n    paragraphs.items[i].style = \"Normal\";\n  }\n}\nawait context.sync();\n\n// --- 2. Swap the in-text citations ---------------------------------------\nasync function replaceAllOccurrences(searchText, replacement) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceAllOccurrences(\"(Karim and Oyewande)\", \"(Lee 208)\");\nawait replaceAllOccurrences(\"(Huang)\", \"(Lee 208)\");\nawait replaceAllOccurrences(\"(Jiang and Ngien)\", \"(Ref-f422876)\");\nawait replaceAllOccurrences(\"(Lee)\", \"(Ref-f422876)\");\nawait replaceAllOccurrences(\"(Naslund et al.)\", \"(Brown & Garcia, 2018)\");\n", "ps1": "# The edit:\n#  1. Three Heading 2 paragraphs (\"Introduction\", \"Immediate Effects of\n#     Social Media\", \"Conclusion\") lose their heading style and revert to\n#     the document's default (\"Normal\") paragraph style.\n#  2. Four in-text citations are swapped out for new ones (a citation\n#     check / reference update):\n#       (Karim and Oyewande)  -> (Lee 208)\n#       (Huang)                -> (Lee 208)\n#       (Jiang and Ngien)      -> (Ref-f422876)\n#       (Lee)                  -> (Ref-f422876)\n#       (Naslund et al.)       -> (Brown & Garcia, 2018)   [occurs twice]\n\n$d = $word.ActiveDocument\n\n# --- 1. Demote the three Heading 2 paragraphs back to Normal -------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Heading 2\") {\n        $p.Style = \"Normal\"\n    }\n}\n\n# --- 2. Swap the in-text citations ---------------------------------------\nfunction Replace-AllText($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nReplace-AllText \"(Karim and Oyewande)\" \"(Lee 208)\"\nReplace-AllText \"(Huang)\" \"(Lee 208)\"\nReplace-AllText \"(Jiang and Ngien)\" \"(Ref-f422876)\"\nReplace-AllText \"(Lee)\" \"(Ref-f422876)\"\nReplace-AllText \"(Naslund et al.)\" \"(Brown & Garcia, 2018)\"\n"}
